$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 10855
$ws.Range("E2").Value = 2158
$ws.Range("F2").Value = 2158
$ws.Range("G2").Value = 2212
$ws.Range("H2").Value = 1683
$ws.Range("I2").Value = 1683
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 9255
$ws.Range("L2").Value = 2275
$ws.Range("M2").Value = 6980
$ws.Range("N2").Value = 6980
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 90
$ws.Range("Q2").Value = 1845
$ws.Range("R2").Value = -699
$ws.Range("S2").Value = -684
$ws.Range("T2").Value = 932
$ws.Range("U2").Value = 913
$ws.Range("V2").ClearContents()
$ws.Range("W2").Value = 19.88
$ws.Range("X2").Value = 15.5
$ws.Range("Y2").Value = 25.5
$ws.Range("Z2").Value = 19.05
$ws.Range("AA2").Value = 32.59
$ws.Range("AB2").Value = 7655.93
$ws.Range("AC2").Value = 3740
$ws.Range("AD2").Value = 21.31
$ws.Range("AE2").Value = 15512
$ws.Range("AF2").Value = 5.14
$ws.Range("AG2").Value = 1670
$ws.Range("AH2").Value = 2.1
$ws.Range("AI2").Value = 44.65
$ws.Range("AJ2").Value = 45000000

# Row 3
$ws.Range("D3").Value = 11797
$ws.Range("E3").Value = 1752
$ws.Range("F3").Value = 1752
$ws.Range("G3").Value = 2213
$ws.Range("H3").Value = 1699
$ws.Range("I3").Value = 1699
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10336
$ws.Range("L3").Value = 2473
$ws.Range("M3").Value = 7863
$ws.Range("N3").Value = 7863
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 90
$ws.Range("Q3").Value = 1398
$ws.Range("R3").Value = -1046
$ws.Range("S3").Value = -752
$ws.Range("T3").Value = 592
$ws.Range("U3").Value = 806
$ws.Range("V3").ClearContents()
$ws.Range("W3").Value = 14.85
$ws.Range("X3").Value = 14.41
$ws.Range("Y3").Value = 22.9
$ws.Range("Z3").Value = 17.35
$ws.Range("AA3").Value = 31.45
$ws.Range("AB3").Value = 8635.879999999999
$ws.Range("AC3").Value = 3776
$ws.Range("AD3").Value = 23.57
$ws.Range("AE3").Value = 17473
$ws.Range("AF3").Value = 5.09
$ws.Range("AG3").Value = 1690
$ws.Range("AH3").Value = 1.9
$ws.Range("AI3").Value = 44.75
$ws.Range("AJ3").Value = 45000000

# Row 4
$ws.Range("D4").Value = 12231
$ws.Range("E4").Value = 1058
$ws.Range("F4").Value = 1058
$ws.Range("G4").Value = 1129
$ws.Range("H4").Value = 883
$ws.Range("I4").Value = 883
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 10925
$ws.Range("L4").Value = 3027
$ws.Range("M4").Value = 7898
$ws.Range("N4").Value = 7898
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 90
$ws.Range("Q4").Value = 1234
$ws.Range("R4").Value = 798
$ws.Range("S4").Value = -760
$ws.Range("T4").Value = 584
$ws.Range("U4").Value = 650
$ws.Range("V4").ClearContents()
$ws.Range("W4").Value = 8.65
$ws.Range("X4").Value = 7.22
$ws.Range("Y4").Value = 11.2
$ws.Range("Z4").Value = 8.300000000000001
$ws.Range("AA4").Value = 38.33
$ws.Range("AB4").Value = 8675.73
$ws.Range("AC4").Value = 1962
$ws.Range("AD4").Value = 27.63
$ws.Range("AE4").Value = 17550
$ws.Range("AF4").Value = 3.09
$ws.Range("AG4").Value = 680
$ws.Range("AH4").Value = 1.25
$ws.Range("AI4").Value = 34.66
$ws.Range("AJ4").Value = 45000000

# Row 5
$ws.Range("D5").Value = 12368
$ws.Range("E5").Value = 1641
$ws.Range("F5").Value = 1641
$ws.Range("G5").Value = 1752
$ws.Range("H5").Value = 1360
$ws.Range("I5").Value = 1360
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 2949
$ws.Range("M5").Value = 9051
$ws.Range("N5").Value = 9051
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 90
$ws.Range("Q5").Value = 1328
$ws.Range("R5").Value = -1934
$ws.Range("S5").Value = -306
$ws.Range("T5").Value = 722
$ws.Range("U5").Value = 607
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 13.26
$ws.Range("X5").Value = 10.99
$ws.Range("Y5").Value = 16.04
$ws.Range("Z5").Value = 11.86
$ws.Range("AA5").Value = 32.59
$ws.Range("AB5").Value = 9963.98
$ws.Range("AC5").Value = 3021
$ws.Range("AD5").Value = 13.44
$ws.Range("AE5").Value = 20113
$ws.Range("AF5").Value = 2.02
$ws.Range("AG5").Value = 1470
$ws.Range("AH5").Value = 3.62
$ws.Range("AI5").Value = 48.65
$ws.Range("AJ5").Value = 45000000

# Row 6
$ws.Range("D6").Value = 12425
$ws.Range("E6").Value = 1915
$ws.Range("F6").Value = 1915
$ws.Range("G6").Value = 2105
$ws.Range("H6").Value = 1613
$ws.Range("I6").Value = 1613
$ws.Range("K6").Value = 12809
$ws.Range("L6").Value = 2968
$ws.Range("M6").Value = 9841
$ws.Range("N6").Value = 9841
$ws.Range("P6").Value = 90
$ws.Range("Q6").Value = 1301
$ws.Range("R6").Value = -1237
$ws.Range("S6").Value = -662
$ws.Range("T6").Value = 484
$ws.Range("U6").Value = 817
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 15.42
$ws.Range("X6").Value = 12.98
$ws.Range("Y6").Value = 17.08
$ws.Range("Z6").Value = 13.01
$ws.Range("AA6").Value = 30.16
$ws.Range("AB6").Value = 10846.87
$ws.Range("AC6").Value = 3585
$ws.Range("AD6").Value = 9.26
$ws.Range("AE6").Value = 21870
$ws.Range("AF6").Value = 1.52
$ws.Range("AG6").Value = 1790
$ws.Range("AH6").Value = 5.39
$ws.Range("AI6").Value = 49.93
$ws.Range("AJ6").Value = 45000000

# Row 7
$ws.Range("D7").Value = 12566
$ws.Range("E7").Value = 1975
$ws.Range("G7").Value = 2061
$ws.Range("H7").Value = 1582
$ws.Range("I7").Value = 1582
$ws.Range("K7").Value = 13758
$ws.Range("L7").Value = 3235
$ws.Range("M7").Value = 10523
$ws.Range("N7").Value = 10536
$ws.Range("P7").Value = 90
$ws.Range("Q7").Value = 1634
$ws.Range("R7").Value = -189
$ws.Range("S7").Value = -658
$ws.Range("T7").Value = 361
$ws.Range("U7").Value = 1328
$ws.Range("W7").Value = 15.72
$ws.Range("X7").Value = 12.59
$ws.Range("Y7").Value = 15.53
$ws.Range("Z7").Value = 11.91
$ws.Range("AA7").Value = 30.74
$ws.Range("AC7").Value = 3516
$ws.Range("AD7").Value = 10.68
$ws.Range("AE7").Value = 23412
$ws.Range("AF7").Value = 1.6
$ws.Range("AG7").Value = 1752
$ws.Range("AH7").Value = 4.67
$ws.Range("AI7").Value = 49.83

# Row 8
$ws.Range("D8").Value = 12801
$ws.Range("E8").Value = 1877
$ws.Range("G8").Value = 2017
$ws.Range("H8").Value = 1547
$ws.Range("I8").Value = 1547
$ws.Range("K8").Value = 14518
$ws.Range("L8").Value = 3219
$ws.Range("M8").Value = 11298
$ws.Range("N8").Value = 11317
$ws.Range("P8").Value = 90
$ws.Range("Q8").Value = 1899
$ws.Range("R8").Value = -585
$ws.Range("S8").Value = -785
$ws.Range("T8").Value = 469
$ws.Range("U8").Value = 1400
$ws.Range("W8").Value = 14.66
$ws.Range("X8").Value = 12.09
$ws.Range("Y8").Value = 14.16
$ws.Range("Z8").Value = 10.94
$ws.Range("AA8").Value = 28.49
$ws.Range("AC8").Value = 3439
$ws.Range("AD8").Value = 10.92
$ws.Range("AE8").Value = 25150
$ws.Range("AF8").Value = 1.49
$ws.Range("AG8").Value = 1705
$ws.Range("AH8").Value = 4.54
$ws.Range("AI8").Value = 49.57

# Row 9
$ws.Range("D9").Value = 13154
$ws.Range("E9").Value = 1926
$ws.Range("G9").Value = 2062
$ws.Range("H9").Value = 1573
$ws.Range("I9").Value = 1573
$ws.Range("K9").Value = 15236
$ws.Range("L9").Value = 3144
$ws.Range("M9").Value = 12092
$ws.Range("N9").Value = 12103
$ws.Range("P9").Value = 90
$ws.Range("Q9").Value = 1861
$ws.Range("R9").Value = -661
$ws.Range("S9").Value = -836
$ws.Range("T9").Value = 539
$ws.Range("U9").Value = 1258
$ws.Range("W9").Value = 14.64
$ws.Range("X9").Value = 11.96
$ws.Range("Y9").Value = 13.43
$ws.Range("Z9").Value = 10.57
$ws.Range("AA9").Value = 26.01
$ws.Range("AC9").Value = 3495
$ws.Range("AD9").Value = 10.74
$ws.Range("AE9").Value = 26897
$ws.Range("AF9").Value = 1.4
$ws.Range("AG9").Value = 1759
$ws.Range("AH9").Value = 4.68
$ws.Range("AI9").Value = 50.34

Write-Output "Applied all IFRS list corrections"